$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block ---
# Name
$ws.Range("B2").Value = "Test Student2"
# Student Number stays "18-0001" (unchanged value, left as-is)
$ws.Range("B3").Value = "18-0001"
# Curriculum Year text
$ws.Range("E2").Value = "Information Technology"
# Year number
$ws.Range("E3").Value = 2022

# Make C3 left-aligned like C2 (cosmetic cleanup matching the shared merge style)
$ws.Range("C3").HorizontalAlignment = $ws.Range("C2").HorizontalAlignment

# --- Remove the old rows 10-12 (no longer needed subjects) ---
$ws.Rows("10:12").Delete()

# --- Rewrite the remaining subject rows (6-9) with the new curriculum data ---
# Row 6
$ws.Range("A6").Value = 214
$ws.Range("B6").Value = "SEMTR"
$ws.Range("C6").Value = "Seminars and Colloquia"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = 1.25
$ws.Range("G6").Value = "Passed"

# Row 7
$ws.Range("A7").Value = 216
$ws.Range("B7").Value = "ST 2"
$ws.Range("C7").Value = "Special Topics 2 (Programming and Database)"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = " "
$ws.Range("F7").Value = 1.5
$ws.Range("G7").Value = "Passed"

# Row 8
$ws.Range("A8").Value = 217
$ws.Range("B8").Value = "TECHNO"
$ws.Range("C8").Value = "Technopreneurship"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = " "
$ws.Range("F8").Value = 1.25
$ws.Range("G8").Value = "Passed"

# Row 9
$ws.Range("A9").Value = 218
$ws.Range("B9").Value = "TECHNO L"
$ws.Range("C9").Value = "Technopreneurship L"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = " "
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = "Passed"

# --- Selection matches the saved view state ---
$ws.Range("B3:C3").Select()
